$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "About"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

# The original (US-model) notes column stays exactly the same, except the
# header in A6 now clarifies these notes came from the US model.
$ws1.Range("A6").Value = 'Notes(from the US model)'

# New column (J) with HK-specific notes explaining the deviations/choices
# made for the Hong Kong model.
$ws1.Range("J6").Value = 'Notes(For this datasheet for HK)):'
$ws1.Range("J7").Value = '1. From the publication of HK Electricity, we know that the gas turbine generators that it owns'
$ws1.Range("J8").Value = 'can be used to fullfill the peak time demand. According to gas generator''s properties, it '
$ws1.Range("J9").Value = 'can also provide flexibility points for solar PV and wind.'
$ws1.Range("J10").Value = 'Therefore, we mark the boolean value of "natural gas peaker" as 1.'
$ws1.Range("J11").Value = '2. Next, as stated in the EPS documentation for elec sector,'
$ws1.Range("J12").Value = '"Which types of plants count as peakers can vary based on model version (country or region represented), '
$ws1.Range("J13").Value = 'but often natural gas peakers and petroleum-fired power plants will qualify".'
$ws1.Range("J14").Value = 'So, although the information from HK''s elec companies didn''t indicate that the oil plant'
$ws1.Range("J15").Value = 'can also be use as peaker, we mark the "pertoleum" 1, same with the US model.'
$ws1.Range("J16").Value = '3.The CLP company in HK has the permission to use 600 MW capacity of the Guangzhou Pumped Storage Power Station,'
$ws1.Range("J17").Value = 'and from the public information we know it is also used for provide flexibility.'
$ws1.Range("J18").Value = 'But in the framwork of EPS model, the flexibility points from Peaker Capacity and Pumped Hydro seems to be '
$ws1.Range("J19").Value = 'calculated at different part. So I think maybe the "hydro" type in this sheet refers not to the "Pumped Hydro",'
$ws1.Range("J20").Value = 'and the flexibility points from Pumed Hydro will be calculated based on BPHC(BAU Pumped Hydro Capacity).'
$ws1.Range("J21").Value = 'For this reason, I didn''t mark the "hydro" as 1 in neither of the two sheets.'
$ws1.Range("J22").Value = '4. To sum up, the current settings(_v1_190405) are totally the same with those in the EPS US model.'

# The new "Notes (HK)" header is bold, matching the style of the other bold
# section headers in column A.
$ws1.Range("J6").Font.Bold = $true

# Leave the same cell selected as in the authored workbook.
$ws1.Range("E19").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "BPaFF-BITPTaP"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("BPaFF-BITPTaP")

# Drop the three rows that used to mirror "natural gas peaker" / "petroleum"
# into crude oil, heavy/residual fuel oil and municipal solid waste -- those
# plant types were removed from this HK-specific sheet.
$ws2.Rows("15:17").Delete()

# The header cell no longer needs the right-aligned "Boolean" style.
$ws2.Range("B1").ClearFormats()

# Highlight, in red, that "hydro" was deliberately left at 0 (explained in
# the new notes column on the About sheet).
$ws2.Range("B5").Font.Color = 255

# Mark that "petroleum" was set to 1 (same as the US model).
$ws2.Range("B11").Font.ColorIndex = 1

$ws2.Range("C23").Select()

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet 3: "BPaFF-BDTPTPF" (identical edits to sheet 2)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")

$ws3.Rows("15:17").Delete()

$ws3.Range("B1").ClearFormats()

$ws3.Range("B5").Font.Color = RGB(255, 0, 0)

$ws3.Range("B11").Font.ColorIndex = 1

$ws3.Range("A23").Select()

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Re-select the About sheet / its original selection so the workbook opens
# where the author left it.
$ws1.Activate()
$ws1.Range("E19").Select()

Write-Output "edits applied"
